# Applies the scheduled-runner market-data refresh to the Leve profit
# columns (H:N) across all eight job sheets, per the commit diff.
# Each row's H-N values are plain (non-formula) numbers sourced from an
# external price feed, so we just overwrite the specific cells that
# changed. A few rows lose their trailing N (or M) cell entirely in the
# new data (no HQ price available) -- those are cleared instead of set.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1293.48
$ws.Range("I28").Value = 1295.2916
$ws.Range("K28").Value = 1295.2916
$ws.Range("M28").Value = -810.2916
$ws.Range("H33").Value = 1910.52
$ws.Range("I33").Value = 2362
$ws.Range("J33").Value = 749.5714
$ws.Range("K33").Value = 2362
$ws.Range("L33").Value = 749.5714
$ws.Range("M33").Value = -2133
$ws.Range("N33").Value = -1207.5714
$ws.Range("H88").Value = 35966.418
$ws.Range("I88").Value = 971.8570999999999
$ws.Range("J88").Value = 84958.8
$ws.Range("K88").Value = 971.8570999999999
$ws.Range("L88").Value = 84958.8
$ws.Range("M88").Value = -565.8570999999999
$ws.Range("N88").Value = -85770.8
$ws.Range("H91").Value = 35966.418
$ws.Range("I91").Value = 971.8570999999999
$ws.Range("J91").Value = 84958.8
$ws.Range("K91").Value = 971.8570999999999
$ws.Range("L91").Value = 84958.8
$ws.Range("M91").Value = 432.1429000000001
$ws.Range("N91").Value = -87766.8
$ws.Range("H97").Value = 1669.2142
$ws.Range("J97").Value = 1669.2142
$ws.Range("L97").Value = 5007.642599999999
$ws.Range("N97").Value = -5999.642599999999
$ws.Range("H111").Value = 1747.6666
$ws.Range("I111").Value = 1747.6666
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 5242.9998
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -2175.9998
$ws.Range("H138").Value = 4027.2727
$ws.Range("I138").Value = 3582.111
$ws.Range("J138").Value = 4086.1912
$ws.Range("K138").Value = 10746.333
$ws.Range("L138").Value = 12258.5736
$ws.Range("M138").Value = -5606.332999999999
$ws.Range("N138").Value = -22538.5736
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 3166
$ws.Range("I4").Value = 2416.0557
$ws.Range("K4").Value = 2416.0557
$ws.Range("M4").Value = -2300.0557
$ws.Range("H55").Value = 21241
$ws.Range("J55").Value = 39052.5
$ws.Range("L55").Value = 39052.5
$ws.Range("N55").Value = -39682.5
$ws.Range("H80").Value = 59996.332
$ws.Range("J80").Value = 59996.332
$ws.Range("L80").Value = 59996.332
$ws.Range("N80").Value = -61992.332
$ws.Range("H83").Value = 59996.332
$ws.Range("J83").Value = 59996.332
$ws.Range("L83").Value = 179988.996
$ws.Range("N83").Value = -189972.996
$ws.Range("H88").Value = 1564.4
$ws.Range("I88").Value = 1391.2
$ws.Range("J88").Value = 1737.6
$ws.Range("K88").Value = 1391.2
$ws.Range("L88").Value = 1737.6
$ws.Range("M88").Value = -985.2
$ws.Range("N88").Value = -2549.6
$ws.Range("H91").Value = 1564.4
$ws.Range("I91").Value = 1391.2
$ws.Range("J91").Value = 1737.6
$ws.Range("K91").Value = 1391.2
$ws.Range("L91").Value = 1737.6
$ws.Range("M91").Value = 12.79999999999995
$ws.Range("N91").Value = -4545.6
$ws.Range("H102").Value = 418590.88
$ws.Range("I102").Value = 501963.5
$ws.Range("K102").Value = 501963.5
$ws.Range("M102").Value = -500341.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4437.7144
$ws.Range("I86").Value = 4734.6665
$ws.Range("J86").Value = 4215
$ws.Range("K86").Value = 4734.6665
$ws.Range("L86").Value = 4215
$ws.Range("M86").Value = -3611.6665
$ws.Range("N86").Value = -6461
$ws.Range("H89").Value = 4437.7144
$ws.Range("I89").Value = 4734.6665
$ws.Range("J89").Value = 4215
$ws.Range("K89").Value = 23673.3325
$ws.Range("L89").Value = 21075
$ws.Range("M89").Value = -18057.3325
$ws.Range("N89").Value = -32307
$ws.Range("H94").Value = 2865.375
$ws.Range("I94").Value = 2218.923
$ws.Range("K94").Value = 2218.923
$ws.Range("M94").Value = -1767.923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4957.375
$ws.Range("I16").Value = 4897.5
$ws.Range("K16").Value = 4897.5
$ws.Range("M16").Value = -4610.5
$ws.Range("H22").Value = 1524.25
$ws.Range("I22").Value = 1199.3334
$ws.Range("K22").Value = 1199.3334
$ws.Range("M22").Value = -849.3334
$ws.Range("H113").Value = 4957.375
$ws.Range("I113").Value = 4897.5
$ws.Range("K113").Value = 4897.5
$ws.Range("M113").Value = -2727.5
$ws.Range("H131").Value = 89053.664
$ws.Range("I131").Value = 83749
$ws.Range("K131").Value = 83749
$ws.Range("M131").Value = -78709

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 792
$ws.Range("I5").Value = 542
$ws.Range("J5").Value = 1208.6666
$ws.Range("K5").Value = 1626
$ws.Range("L5").Value = 3625.9998
$ws.Range("M5").Value = -1514
$ws.Range("N5").Value = -3849.9998
$ws.Range("H132").Value = 3703.5
$ws.Range("I132").Value = 1203.6
$ws.Range("J132").Value = 5489.143
$ws.Range("K132").Value = 10832.4
$ws.Range("L132").Value = 49402.287
$ws.Range("M132").Value = -8302.4
$ws.Range("N132").Value = -54462.287
$ws.Range("H135").Value = 792
$ws.Range("I135").Value = 542
$ws.Range("J135").Value = 1208.6666
$ws.Range("K135").Value = 4878
$ws.Range("L135").Value = 10877.9994
$ws.Range("M135").Value = -2343
$ws.Range("N135").Value = -15947.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9937.317999999999
$ws.Range("J80").Value = 13352.066
$ws.Range("L80").Value = 13352.066
$ws.Range("N80").Value = -15348.066
$ws.Range("H83").Value = 9937.317999999999
$ws.Range("J83").Value = 13352.066
$ws.Range("L83").Value = 66760.33
$ws.Range("N83").Value = -76744.33
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H132").Value = 5953.5264
$ws.Range("I132").Value = 4017.1667
$ws.Range("J132").Value = 9273
$ws.Range("K132").Value = 12051.5001
$ws.Range("L132").Value = 27819
$ws.Range("M132").Value = -9521.500100000001
$ws.Range("N132").Value = -32879
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4529.242
$ws.Range("I122").Value = 4170.0415
$ws.Range("J122").Value = 5487.1113
$ws.Range("K122").Value = 12510.1245
$ws.Range("L122").Value = 16461.3339
$ws.Range("M122").Value = -10060.1245
$ws.Range("N122").Value = -21361.3339
$ws.Range("H136").Value = 3335.0833
$ws.Range("I136").Value = 2513.2954
$ws.Range("K136").Value = 7539.8862
$ws.Range("M136").Value = -4989.8862

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 21995
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H81").Value = 4927.4287
$ws.Range("I81").Value = 1831
$ws.Range("J81").Value = 7249.75
$ws.Range("K81").Value = 3662
$ws.Range("L81").Value = 14499.5
$ws.Range("M81").Value = -2601
$ws.Range("N81").Value = -16621.5
$ws.Range("H84").Value = 4927.4287
$ws.Range("I84").Value = 1831
$ws.Range("J84").Value = 7249.75
$ws.Range("K84").Value = 18310
$ws.Range("L84").Value = 72497.5
$ws.Range("M84").Value = -13006
$ws.Range("N84").Value = -83105.5
$ws.Range("H100").Value = 1833418.6
$ws.Range("I100").Value = 2748221
$ws.Range("K100").Value = 5496442
$ws.Range("M100").Value = -5495901
$ws.Range("N31").ClearContents()

